# PPT => NG 8
# Slide 5, content placeholder: change "ng new angular-tour-of-heroes"
# to "ng new toh" by merging the "angular" / "-tour-of-" / "heroes" runs
# into a single run with text "toh" (keeping the first run's formatting).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 4 is "ng new angular-tour-of-heroes"
$para = $tr.Paragraphs(4, 1)

# Run 3 is "angular" -> becomes "toh"
$runAngular = $para.Runs(3, 1)
$runAngular.Text = "toh"

# What is now run 4 is the old "-tour-of-" run -> remove it
$runTourOf = $para.Runs(4, 1)
$runTourOf.Text = ""

# What is now run 4 is the old "heroes" run -> remove it
$runHeroes = $para.Runs(4, 1)
$runHeroes.Text = ""
